$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19: resale numbers update for 2024-01-04 22:53:35.
# Columns A-D hold text-like values ("2024-01-04", "22:53:35", "Thursday",
# "00"); without forcing Text format Excel would auto-coerce these into a
# date serial / time / number. Format as Text first, then clear the
# resulting formatting so the cells end up with no explicit style (matching
# the plain data rows above), while the stored values remain text.
$textRange = $ws.Range("A19:D19")
$textRange.NumberFormat = "@"

$ws.Range("A19").Value = "2024-01-04"
$ws.Range("B19").Value = "22:53:35"
$ws.Range("C19").Value = "Thursday"
$ws.Range("D19").Value = "00"

$textRange.ClearFormats()

$ws.Range("E19").Value = 140535
$ws.Range("F19").Value = 142882
$ws.Range("G19").Value = 172309
$ws.Range("H19").Value = 147196
$ws.Range("I19").Value = -1
$ws.Range("J19").Value = 118040
$ws.Range("K19").Value = 224369
$ws.Range("L19").Value = 248599
$ws.Range("M19").Value = 184650
$ws.Range("N19").Value = 110083
$ws.Range("O19").Value = 40410
$ws.Range("P19").Value = 30798
$ws.Range("Q19").Value = 72395
$ws.Range("R19").Value = -1
$ws.Range("S19").Value = 41789
$ws.Range("T19").Value = -1
